$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-24 18:25:38"

# Update the timestamp in column A for existing rows 2-12
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Row 13: new entry
$ws.Cells.Item(13, 1).Value = $newTimestamp
$ws.Cells.Item(13, 2).Value = "【急募】イベント用問い合わせLINE構築のフリーランス募集!"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5420186"
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://www.lancers.jp/work/detail/5420186")
$ws.Cells.Item(13, 6).Style = $ws.Cells.Item(12, 6).Style
$ws.Cells.Item(13, 7).Value = 10

# Row 14: new entry
$ws.Cells.Item(14, 1).Value = $newTimestamp
$ws.Cells.Item(14, 2).Value = "【急募】Google Play Consoleでのクローズテスト実施者募集!"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5419425"
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://www.lancers.jp/work/detail/5419425")
$ws.Cells.Item(14, 6).Style = $ws.Cells.Item(12, 6).Style
$ws.Cells.Item(14, 7).Value = 10
